$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (old 2:21) down to (3:22)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with its data
$ws.Range("A2").Value2 = -0.1041525229811668
$ws.Range("B2").Value2 = 0.2035708427429199
$ws.Range("C2").Value2 = 0.0946841165423393

# Append 9 new data rows after the (now shifted) last row, i.e. rows 23-31
$newRows = @(
    @(-0.6001750826835632, -0.166460782289505,  -0.2874121069908142),
    @(-0.1484402567148208,  0.0137444678694009,  0.1325577646493911),
    @(-0.0577267669141292, -0.1240056455135345,  0.0187841057777404),
    @( 0.0021380283869802, -0.0429132841527462,  0.1725694388151168),
    @( 0.0100792767480015,  0.0042760567739605,  0.0404698215425014),
    @(-0.0004581489483825, -0.0215329993516206, -0.0074830991216003),
    @(-0.0334448739886283, -0.1394299864768982, -0.0435241498053073),
    @(-0.040775254368782,  -0.2396118938922882, -0.1838704347610473),
    @( 0.057115901261568,  -0.2092213481664657, -0.0080939643085002)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 2).Value2 = $vals[1]
    $ws.Cells.Item($row, 3).Value2 = $vals[2]
}
